$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# --- Criterion 1 Online collaboration block (rows 2-11) ---

# Row 5: was placeholder "Peer assessment 3" -> becomes Alex's peer assessment
$ws.Cells.Item(5, 1).Value = "Alex"
$ws.Cells.Item(5, 2).Value = "Good"
$ws.Cells.Item(5, 3).Value = "Research, Hardware setup "

# Row 2: Self assesment row - fill in grade + comment
$ws.Cells.Item(2, 2).Value = "Good"
$ws.Cells.Item(2, 3).Value = @"
I have been active in the weekly online Discord meetings - both in 
terms of planning (Adding discussion points to the meeting agenda, taking responsibility for leading a group discussion or a walkthrough of conducted work, etc) 
I have been reading through literature and materials uploaded
to GiT by other project group members, and actively been adding comments and suggestions. I have been uploading relevant materials to
GiT as well.  
I have been responding fairly quickly to messages, both private and group announcements, on the discord platform. 
"@

# --- Criterion 1 International Collaboration block (rows 62-71) ---

# Row 63: was placeholder "Peer assessment 1" -> becomes Alex's peer assessment
$ws.Cells.Item(63, 1).Value = "Alex"
$ws.Cells.Item(63, 2).Value = "Excellent"
$ws.Cells.Item(63, 3).Value = "Active collaborator, Motivated "

# Row 62: Self assesment row - fill in grade only
$ws.Cells.Item(62, 2).Value = "Good"

# Match the author's final selection/scroll position in the sheet view
$ws.Range("B63").Select()
